# Applies the "all canine test cases 72" commit:
#   - Adds three new worksheets: CypherOutput_Message, StatOutput, StatOutput_Message
#   - CypherOutput_Message is a duplicate of the existing "Message" worksheet
#   - StatOutput holds a small 2x4 results table (counts for the new "stat" cypher query)
#   - StatOutput_Message is the "Message" worksheet content written twice (rows 1-10 and
#     11-20), with the second copy's Cypher-text row replaced by the new stat query text

$wb = $excel.ActiveWorkbook

$msgSheet = $wb.Worksheets.Item("Message")

# --- New stat query text (used later, in StatOutput_Message row 18) ---
$statCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN ['NCATS-COP01']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# ---------------------------------------------------------------------------
# 1) CypherOutput_Message -- exact duplicate of the "Message" sheet, placed
#    right after it.
# ---------------------------------------------------------------------------
$msgSheet.Copy($null, $msgSheet)
$cypherMsgSheet = $wb.Worksheets.Item($msgSheet.Index + 1)
$cypherMsgSheet.Name = "CypherOutput_Message"

# ---------------------------------------------------------------------------
# 2) StatOutput -- small 2x4 table of counts, placed after CypherOutput_Message.
# ---------------------------------------------------------------------------
$statSheet = $wb.Worksheets.Add($null, $cypherMsgSheet)
$statSheet.Name = "StatOutput"

$statSheet.Range("A1").Value = "number_of_files"
$statSheet.Range("B1").Value = "number_of_sample"
$statSheet.Range("C1").Value = "number_of_cases"
$statSheet.Range("D1").Value = "number_of_study"

# These look numeric but must be written as text (matches source export which
# always emits plain <t> shared-string cells for its results table).
$statValuesRange = $statSheet.Range("A2:D2")
$statValuesRange.NumberFormat = "@"
$statSheet.Range("A2").Value = "331"
$statSheet.Range("B2").Value = "136"
$statSheet.Range("C2").Value = "60"
$statSheet.Range("D2").Value = "1"

# ---------------------------------------------------------------------------
# 3) StatOutput_Message -- the "Message" sheet content written twice
#    (rows 1-10 then rows 11-20), with the second copy's Cypher-text cell
#    (row 18, the 8th row of the second block) swapped for the new stat query.
# ---------------------------------------------------------------------------
$statMsgSheet = $wb.Worksheets.Add($null, $statSheet)
$statMsgSheet.Name = "StatOutput_Message"

$msgBlock = $msgSheet.Range("A1:A10")
$msgBlock.Copy()
$statMsgSheet.Range("A1").PasteSpecial()
$msgBlock.Copy()
$statMsgSheet.Range("A11").PasteSpecial()

$statMsgSheet.Range("A18").Value = $statCypher
